# Update the computed (output) values on sheet "A-09", row 2, reflecting
# the results of running cases A-F for the first time. Row 2 holds a single
# simulation case; the input cells (round numbers such as 40, 65, 50, 10,
# -10) stay the same while all derived/computed quantities are refreshed
# with newly calculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A-09")

$ws.Range("B2").Value = -9424.810314801476
$ws.Range("C2").Value = 13561.10628077167
$ws.Range("D2").Value = -3520.258773335945
$ws.Range("E2").Value = -616.0371926342616
$ws.Range("F2").Value = 55.96194096571918
$ws.Range("G2").Value = 40.92402927951355
$ws.Range("H2").Value = 56.24114139409272
$ws.Range("I2").Value = 45.42231789697769
$ws.Range("J2").Value = 55.96194096582354
$ws.Range("K2").Value = 50.34513396833665
$ws.Range("L2").Value = 53.46401473543068
$ws.Range("O2").Value = 45.62699102746905
$ws.Range("P2").Value = 55.02894110476564
$ws.Range("R2").Value = 8.953380639800754
$ws.Range("S2").Value = -17.90676128295929
$ws.Range("T2").Value = 8.953380643158548
$ws.Range("X2").Value = -92.33646748703023
$ws.Range("Y2").Value = -142.3360080380841
$ws.Range("Z2").Value = -92.33646748703015
$ws.Range("AE2").Value = -8.953380639800754
$ws.Range("AF2").Value = 8.953380643158548
$ws.Range("AG2").Value = 8.953380639800754
$ws.Range("AH2").Value = -17.90676128295929
$ws.Range("AI2").Value = 8.953380643158548
$ws.Range("AJ2").Value = 8.953380639800754
$ws.Range("AK2").Value = -8.953380643158548
$ws.Range("AL2").Value = 24.99977027552693
$ws.Range("AM2").Value = -24.99977027552692
$ws.Range("AN2").Value = -92.33646748703023
$ws.Range("AO2").Value = -142.3360080380841
$ws.Range("AP2").Value = -92.33646748703015
$ws.Range("AQ2").Value = -24.99977027552693
$ws.Range("AR2").Value = 24.99977027552692
$ws.Range("AS2").Value = 55.96194096571918
$ws.Range("AT2").Value = 55.96194096571918
$ws.Range("AU2").Value = 56.24114139409255
$ws.Range("AV2").Value = 56.24114139409255
$ws.Range("AW2").Value = 56.24114139409272
$ws.Range("AX2").Value = 55.96194096582354
$ws.Range("AY2").Value = 55.96194096582354
$ws.Range("AZ2").Value = 40.92402927951355
$ws.Range("BA2").Value = 40.92402927951355
$ws.Range("BB2").Value = 45.42231789697769
$ws.Range("BC2").Value = 40.74149327093306
$ws.Range("BD2").Value = 50.10314252126648
$ws.Range("BE2").Value = 50.34513396833665
$ws.Range("BF2").Value = 50.34513396833665
$ws.Range("BG2").Value = 53.46401473543068
$ws.Range("BJ2").Value = 45.62699102746905
$ws.Range("BK2").Value = 55.02894110476564
